$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet (sheet1): insert a new row for "2022-Q3" right under the
#    header row, push the existing quarters down, and renumber the index
#    column (A) so it stays 0,1,2,3.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows(2).Insert()

# Copy the formatting of the row that just got pushed down to row 3 (the
# old "2022-Q2" row) into the brand new row 2, so the new row looks just
# like the others (bordered/bold index cell in column A, plain data cells).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 3
$summary.Cells.Item(2, 4).Value = 0.15

# Renumber the index column for the rows that shifted down.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3

# ---------------------------------------------------------------------------
# 2. Add the new "2022-Q3" worksheet right after "总计", using the old
#    "2021-Q4" sheet (5 data rows) as a formatting template since it has
#    enough rows to cover the 3 data rows the new sheet needs, then trim the
#    extra rows and overwrite the values.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(4)
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template had 5 data rows (2..6); we only need 3 (2..4).
$q3.Rows("5:6").Delete()

function Set-TextCell($sheet, $addr, $text) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $text
}

# Row 2 -> 233009 / 大摩多因子精选策略混合
$q3.Cells.Item(2, 1).Value = 0
Set-TextCell $q3 "B2" "233009"
$q3.Cells.Item(2, 3).Value = "大摩多因子精选策略混合"
Set-TextCell $q3 "D2" "6.50"
Set-TextCell $q3 "E2" "83.44"
Set-TextCell $q3 "F2" "0.93"
Set-TextCell $q3 "G2" "0.0604"
$q3.Cells.Item(2, 8).Value = 7

# Row 3 -> 516620 / 国泰中证影视主题ETF
$q3.Cells.Item(3, 1).Value = 1
Set-TextCell $q3 "B3" "516620"
$q3.Cells.Item(3, 3).Value = "国泰中证影视主题ETF"
Set-TextCell $q3 "D3" "0.94"
Set-TextCell $q3 "E3" "99.07"
Set-TextCell $q3 "F3" "4.83"
Set-TextCell $q3 "G3" "0.0454"
$q3.Cells.Item(3, 8).Value = 8

# Row 4 -> 159855 / 银华中证影视主题ETF
$q3.Cells.Item(4, 1).Value = 2
Set-TextCell $q3 "B4" "159855"
$q3.Cells.Item(4, 3).Value = "银华中证影视主题ETF"
Set-TextCell $q3 "D4" "0.84"
Set-TextCell $q3 "E4" "96.84"
Set-TextCell $q3 "F4" "4.70"
Set-TextCell $q3 "G4" "0.0395"
$q3.Cells.Item(4, 8).Value = 8
